$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill H1305:I1571 with the value 999, matching the target edit.
$ws.Range("H1305:I1571").Value = 999
